$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 406.83334
$ws.Range("I2").Value = 410.17648
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 410.17648
$ws.Range("L2").Value = 350
$ws.Range("M2").Value = -297.17648
$ws.Range("N2").Value = -576

$ws.Range("H4").Value = 730.3333
$ws.Range("I4").Value = 647.5
$ws.Range("J4").Value = 896
$ws.Range("K4").Value = 647.5
$ws.Range("L4").Value = 896
$ws.Range("M4").Value = -533.5
$ws.Range("N4").Value = -1124

$ws.Range("H112").Value = 2011.8572
$ws.Range("J112").Value = 2107.842
$ws.Range("L112").Value = 6323.526
$ws.Range("N112").Value = -8539.526

$ws.Range("H125").Value = 2631
$ws.Range("I125").Value = 2642
$ws.Range("J125").Value = 2622.3572
$ws.Range("K125").Value = 23778
$ws.Range("L125").Value = 23601.2148
$ws.Range("M125").Value = -21318
$ws.Range("N125").Value = -28521.2148

$ws.Range("H138").Value = 2740.883
$ws.Range("J138").Value = 2762.5293
$ws.Range("L138").Value = 8287.5879
$ws.Range("N138").Value = -18567.5879

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2449.9285
$ws.Range("I2").Value = 2483.1667
$ws.Range("J2").Value = 2425
$ws.Range("K2").Value = 2483.1667
$ws.Range("L2").Value = 2425
$ws.Range("M2").Value = -2370.1667
$ws.Range("N2").Value = -2651

$ws.Range("H116").Value = 2449.9285
$ws.Range("I116").Value = 2483.1667
$ws.Range("J116").Value = 2425
$ws.Range("K116").Value = 2483.1667
$ws.Range("L116").Value = 2425
$ws.Range("M116").Value = -189.1667000000002
$ws.Range("N116").Value = -7013

$ws.Range("H132").Value = 2829.1904
$ws.Range("I132").Value = 2474.2666
$ws.Range("J132").Value = 3716.5
$ws.Range("K132").Value = 7422.7998
$ws.Range("L132").Value = 11149.5
$ws.Range("M132").Value = -4892.7998
$ws.Range("N132").Value = -16209.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2449.9285
$ws.Range("I3").Value = 2483.1667
$ws.Range("J3").Value = 2425
$ws.Range("K3").Value = 2483.1667
$ws.Range("L3").Value = 2425
$ws.Range("M3").Value = -2369.1667
$ws.Range("N3").Value = -2653

$ws.Range("H86").Value = 38464412
$ws.Range("I86").Value = 45457176
$ws.Range("J86").Value = 4226.75
$ws.Range("K86").Value = 45457176
$ws.Range("L86").Value = 4226.75
$ws.Range("M86").Value = -45456053
$ws.Range("N86").Value = -6472.75

$ws.Range("H89").Value = 38464412
$ws.Range("I89").Value = 45457176
$ws.Range("J89").Value = 4226.75
$ws.Range("K89").Value = 227285880
$ws.Range("L89").Value = 21133.75
$ws.Range("M89").Value = -227280264
$ws.Range("N89").Value = -32365.75

$ws.Range("H94").Value = 15625997
$ws.Range("I94").Value = 19231726
$ws.Range("J94").Value = 1166.6666
$ws.Range("K94").Value = 19231726
$ws.Range("L94").Value = 1166.6666
$ws.Range("M94").Value = -19231275
$ws.Range("N94").Value = -2068.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1705.5593
$ws.Range("I31").Value = 1668.2545
$ws.Range("K31").Value = 1668.2545
$ws.Range("M31").Value = -1373.2545

$ws.Range("H34").Value = 1705.5593
$ws.Range("I34").Value = 1668.2545
$ws.Range("K34").Value = 1668.2545
$ws.Range("M34").Value = -1466.2545

$ws.Range("H94").Value = 1881.7059
$ws.Range("I94").Value = 1408.75
$ws.Range("J94").Value = 2302.111
$ws.Range("K94").Value = 1408.75
$ws.Range("L94").Value = 2302.111
$ws.Range("M94").Value = -957.75
$ws.Range("N94").Value = -3204.111

$ws.Range("H105").Value = 1300
$ws.Range("I105").Value = 1800
$ws.Range("J105").Value = 800
$ws.Range("K105").Value = 1800
$ws.Range("L105").Value = 800
$ws.Range("M105").Value = -53
$ws.Range("N105").Value = -4294

$ws.Range("H141").Value = 340162.9
$ws.Range("J141").Value = 340162.9
$ws.Range("L141").Value = 340162.9
$ws.Range("N141").Value = -350522.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 253.26315
$ws.Range("I14").Value = 253.26315
$ws.Range("K14").Value = 759.78945
$ws.Range("M14").Value = -586.78945

$ws.Range("H34").Value = 8334659
$ws.Range("I34").Value = 151.16667
$ws.Range("J34").Value = 16669167
$ws.Range("K34").Value = 453.50001
$ws.Range("L34").Value = 50007501
$ws.Range("M34").Value = -369.50001
$ws.Range("N34").Value = -50007669

$ws.Range("H39").Value = 4087.76
$ws.Range("J39").Value = 4112.7827
$ws.Range("L39").Value = 12338.3481
$ws.Range("N39").Value = -12926.3481

$ws.Range("H55").Value = 2210
$ws.Range("J55").Value = 2857.1428
$ws.Range("L55").Value = 8571.428400000001
$ws.Range("N55").Value = -8925.428400000001

$ws.Range("H131").Value = 23846218
$ws.Range("I131").Value = 111111550
$ws.Range("J131").Value = 46580.88
$ws.Range("K131").Value = 333334650
$ws.Range("L131").Value = 139742.64
$ws.Range("M131").Value = -333329610
$ws.Range("N131").Value = -149822.64

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7159.8335
$ws.Range("I132").Value = 9060.666999999999
$ws.Range("K132").Value = 27182.001
$ws.Range("M132").Value = -24652.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1499.75
$ws.Range("I16").Value = 1499.75
$ws.Range("K16").Value = 1499.75
$ws.Range("M16").Value = -1329.75

$ws.Range("H132").Value = 2491.1292
$ws.Range("I132").Value = 2054.2354
$ws.Range("J132").Value = 3021.6428
$ws.Range("K132").Value = 6162.706200000001
$ws.Range("L132").Value = 9064.928400000001
$ws.Range("M132").Value = -3632.706200000001
$ws.Range("N132").Value = -14124.9284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 32010
$ws.Range("J118").Value = 32010
$ws.Range("L118").Value = 32010
$ws.Range("N118").Value = -35324

$ws.Range("H132").Value = 3192.58
$ws.Range("I132").Value = 3448.0789
$ws.Range("J132").Value = 2383.5
$ws.Range("K132").Value = 10344.2367
$ws.Range("L132").Value = 7150.5
$ws.Range("M132").Value = -7814.236699999999
$ws.Range("N132").Value = -12210.5

$ws.Range("H133").Value = 26990.834
$ws.Range("J133").Value = 26990.834
$ws.Range("L133").Value = 26990.834
$ws.Range("N133").Value = -37110.834
